$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels: "_old" -> "_FV2210", "_new" -> "_FV2304" ------
$headerMap = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}
foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2. Freeze the header row (pane split after row 1) ----------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn A1:U82 into a native Excel Table named "Table1" ---------------
# The header row (A1:U1) already carries explicit formatting (bold font,
# fill, border -- style index "s=1"). ListObjects.Add() auto-generates a
# header dxf (bold) whenever the header range carries ANY non-default
# style at the moment of creation, which would add an unwanted <dxfs>
# entry. Work around it by: stashing the header's format on a scratch
# cell, clearing the header's format (now default / dxf-free), creating
# the table, then pasting the stashed format back and discarding the
# scratch cell -- leaving styles.xml untouched.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("W1")
$ws.Range("A1").Copy($scratch)
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U82")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

$scratch.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$scratch.Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

Write-Host "Edit applied"
